$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.239067055393586
$ws.Range("C2").Value = 0.4518950437317784
$ws.Range("J2").Value = 0.008746355685131196
$ws.Range("P2").Value = 0.1778425655976676
$ws.Range("S2").Value = 0.1224489795918367
$ws.Range("B3").Value = 0.01886792452830189
$ws.Range("C3").Value = 0.02515723270440252
$ws.Range("J3").Value = 0.006289308176100629
$ws.Range("P3").Value = 0.7358490566037735
$ws.Range("S3").Value = 0.2138364779874214
$ws.Range("J4").Value = 0.04444444444444445
$ws.Range("P4").Value = 0.6222222222222222
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.09313725490196079
$ws.Range("F6").Value = 0.0392156862745098
$ws.Range("J6").Value = 0.2745098039215687
$ws.Range("O6").Value = 0.03431372549019608
$ws.Range("Q6").Value = 0.2009803921568628
$ws.Range("R6").Value = 0.05392156862745098
$ws.Range("S6").Value = 0.303921568627451
$ws.Range("B7").Value = 0.1262135922330097
$ws.Range("D7").Value = 0.02912621359223301
$ws.Range("F7").Value = 0.07281553398058252
$ws.Range("J7").Value = 0.1359223300970874
$ws.Range("O7").Value = 0.02427184466019417
$ws.Range("Q7").Value = 0.2038834951456311
$ws.Range("R7").Value = 0.06796116504854369
$ws.Range("S7").Value = 0.3398058252427185
$ws.Range("B8").Value = 0.122969837587007
$ws.Range("D8").Value = 0.01624129930394431
$ws.Range("E8").Value = 0.002320185614849188
$ws.Range("F8").Value = 0.04872389791183294
$ws.Range("J8").Value = 0.1345707656612529
$ws.Range("O8").Value = 0.01624129930394431
$ws.Range("Q8").Value = 0.185614849187935
$ws.Range("R8").Value = 0.06960556844547564
$ws.Range("S8").Value = 0.4037122969837587
$ws.Range("B9").Value = 0.08530805687203792
$ws.Range("D9").Value = 0.02369668246445497
$ws.Range("F9").Value = 0.06635071090047394
$ws.Range("J9").Value = 0.1327014218009479
$ws.Range("O9").Value = 0.02843601895734597
$ws.Range("Q9").Value = 0.1990521327014218
$ws.Range("R9").Value = 0.07109004739336493
$ws.Range("S9").Value = 0.3933649289099526
$ws.Range("B10").Value = 0.1061410159211524
$ws.Range("D10").Value = 0.021986353297953
$ws.Range("E10").Value = 0.0007581501137225171
$ws.Range("F10").Value = 0.06141015921152388
$ws.Range("J10").Value = 0.1175132676269901
$ws.Range("O10").Value = 0.01213040181956027
$ws.Range("Q10").Value = 0.2524639878695982
$ws.Range("R10").Value = 0.05913570887035633
$ws.Range("S10").Value = 0.3684609552691433
$ws.Range("G11").Value = 0.129746835443038
$ws.Range("J11").Value = 0.1075949367088608
$ws.Range("K11").Value = 0.189873417721519
$ws.Range("L11").Value = 0.5632911392405063
$ws.Range("S11").Value = 0.00949367088607595
$ws.Range("G12").Value = 0.7675675675675676
$ws.Range("J12").Value = 0.1621621621621622
$ws.Range("K12").Value = 0.01081081081081081
$ws.Range("L12").Value = 0.01621621621621622
$ws.Range("S12").Value = 0.04324324324324325
$ws.Range("G13").Value = 0.6410256410256411
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("S13").Value = 0.02564102564102564
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.03940886699507389
$ws.Range("H15").Value = 0.167487684729064
$ws.Range("I15").Value = 0.06896551724137931
$ws.Range("J15").Value = 0.3349753694581281
$ws.Range("K15").Value = 0.0541871921182266
$ws.Range("M15").Value = 0.01970443349753695
$ws.Range("N15").Value = 0.004926108374384237
$ws.Range("O15").Value = 0.06403940886699508
$ws.Range("S15").Value = 0.2463054187192118
$ws.Range("F16").Value = 0.025
$ws.Range("H16").Value = 0.145
$ws.Range("I16").Value = 0.1
$ws.Range("J16").Value = 0.465
$ws.Range("K16").Value = 0.12
$ws.Range("M16").Value = 0.01
$ws.Range("O16").Value = 0.02
$ws.Range("S16").Value = 0.115
$ws.Range("F17").Value = 0.01310861423220974
$ws.Range("H17").Value = 0.1610486891385768
$ws.Range("I17").Value = 0.09737827715355805
$ws.Range("J17").Value = 0.4625468164794008
$ws.Range("K17").Value = 0.08052434456928839
$ws.Range("M17").Value = 0.02059925093632959
$ws.Range("O17").Value = 0.05430711610486891
$ws.Range("S17").Value = 0.1104868913857678
$ws.Range("F18").Value = 0.0272108843537415
$ws.Range("H18").Value = 0.1564625850340136
$ws.Range("I18").Value = 0.1292517006802721
$ws.Range("J18").Value = 0.3469387755102041
$ws.Range("K18").Value = 0.1428571428571428
$ws.Range("M18").Value = 0.0272108843537415
$ws.Range("O18").Value = 0.04761904761904762
$ws.Range("S18").Value = 0.1224489795918367
$ws.Range("F19").Value = 0.01528559935639582
$ws.Range("H19").Value = 0.2123893805309734
$ws.Range("I19").Value = 0.08608205953338696
$ws.Range("J19").Value = 0.3781174577634754
$ws.Range("K19").Value = 0.1238938053097345
$ws.Range("M19").Value = 0.01528559935639582
$ws.Range("N19").Value = 0.0008045052292839903
$ws.Range("O19").Value = 0.06114239742558326
$ws.Range("S19").Value = 0.1069991954947707
